$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.974.51"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "2.298.52"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'317.42"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'104.53"
$ws.Range("E6").Value = "  -1.70%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").Value = "'39.42"
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "'8.38"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").Value = "'0.973"
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").Value = "'15.39"
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("D16").Value = "2.646.65"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").Value = "2.303.70"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "42.076.94"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'7.66"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "'288.78"
$ws.Range("E21").Value = "  +13.05%  "
$ws.Range("D22").Value = "'73.74"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'2.28"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").Value = "'10.03"
$ws.Range("E25").Value = "  +7.21%  "
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").Value = "'10.91"
$ws.Range("E27").Value = "  -3.97%  "
$ws.Range("D28").Value = "'23.43"
$ws.Range("E28").Value = "  +2.97%  "
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").Value = "'164.33"
$ws.Range("E30").Value = "  -6.55%  "
$ws.Range("D31").Value = "'35.26"
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("D32").Value = "'0.0883"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").Value = "'5.87"
$ws.Range("E34").Value = "  -2.94%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("E36").Value = "  -5.84%  "
$ws.Range("D37").Value = "'4.62"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("D38").Value = "'0.0351"
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").Value = "'2.84"
$ws.Range("E39").Value = "  +6.89%  "
$ws.Range("E40").Value = "  -5.82%  "
$ws.Range("D41").Value = "'103.46"
$ws.Range("E41").Value = "  +21.88%  "
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("D43").Value = "'70.75"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'0.226"
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "'116.36"
$ws.Range("E46").Value = "  +3.70%  "
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("D48").Value = "'78.34"
$ws.Range("E48").Value = "  +7.16%  "
$ws.Range("D49").Value = "'9.08"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "'5.33"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("E51").Value = "  +0.09%  "
